$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$changes = @{
    9  = -7
    10 = -1
    11 = -2
    12 = 5
    16 = -5
    17 = 0
    20 = -4
    23 = -3
    24 = -4
    26 = 2
    27 = -5
    28 = 4
    29 = -2
    36 = -9
    37 = -2
    39 = -2
    42 = -3
    43 = 2
    49 = 2
    51 = -3
    55 = -4
    58 = -2
    62 = 2
    66 = 1
    67 = -1
    68 = 7
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
